$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $text) {
    # Force the cell to be treated as text even when the string looks like
    # a pure number (e.g. "305.47"), then restore the default "Normal"
    # style so no stray number-format styling is left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "22.510.72"
$ws.Range("E2").Value = "  +9.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.608.90"
$ws.Range("E3").Value = "  +9.03%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.56%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "305.47"
$ws.Range("E5").Value = "  +9.86%  "

# Row 6 - USDC
Set-TextValue $ws.Range("D6") "0.9919"
$ws.Range("E6").Value = "  +4.45%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.3677"
$ws.Range("E7").Value = "  +1.55%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3402"
$ws.Range("E8").Value = "  +11.09%  "

# Row 9 - OKB
Set-TextValue $ws.Range("D9") "42.30"
$ws.Range("E9").Value = "  +7.23%  "

# Row 10 - Polygon
$ws.Range("E10").Value = "  +7.51%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.07062"
$ws.Range("E11").Value = "  +6.12%  "

# Row 12 - BinanceUSD
Set-TextValue $ws.Range("D12") "0.9997"
$ws.Range("E12").Value = "  -0.25%  "

# Row 13 - Solana
Set-TextValue $ws.Range("D13") "19.74"
$ws.Range("E13").Value = "  +8.92%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "5.914"
$ws.Range("E14").Value = "  +7.07%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "6.639"
$ws.Range("E15").Value = "  +7.12%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "1.610.18"
$ws.Range("E16").Value = "  +9.39%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  +5.37%  "

# Row 18 - Dai
Set-TextValue $ws.Range("D18") "0.9921"
$ws.Range("E18").Value = "  +4.38%  "

# Row 19 - TRON
Set-TextValue $ws.Range("D19") "0.06698"
$ws.Range("E19").Value = "  +12.80%  "

# Row 20 - Litecoin
Set-TextValue $ws.Range("D20") "78.01"
$ws.Range("E20").Value = "  +12.45%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.017"
$ws.Range("E21").Value = "  +9.36%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +10.98%  "

# Row 23 - Cosmos
Set-TextValue $ws.Range("D23") "11.85"
$ws.Range("E23").Value = "  +6.66%  "

# Row 24 - WrappedBTC
$ws.Range("D24").Value = "22.514.44"
$ws.Range("E24").Value = "  +9.34%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.387"
$ws.Range("E25").Value = "  +5.45%  "

# Row 26 - LidoDAOToken
Set-TextValue $ws.Range("D26") "2.592"
$ws.Range("E26").Value = "  +22.03%  "

# Row 27 - Monero
Set-TextValue $ws.Range("D27") "149.43"
$ws.Range("E27").Value = "  +4.66%  "

# Row 28 - EthereumClassic
Set-TextValue $ws.Range("D28") "19.55"
$ws.Range("E28").Value = "  +13.59%  "

# Row 29 - WrappedliquidstakedEther2.0
$ws.Range("D29").Value = "1.792.40"
$ws.Range("E29").Value = "  +9.77%  "

# Row 30 - BitcoinCash
Set-TextValue $ws.Range("D30") "122.73"
$ws.Range("E30").Value = "  +7.96%  "

# Row 31 - was HuobiToken, now Filecoin
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D31") "6.155"
$ws.Range("E31").Value = "  +22.52%  "

# Row 32 - was Filecoin, now HuobiToken
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D32") "4.007"
$ws.Range("E32").Value = "  +1.47%  "

# Row 33 - ImmutableX
Set-TextValue $ws.Range("D33") "0.9523"
$ws.Range("E33").Value = "  +17.51%  "

# Row 34 - WEMIXTOKEN
Set-TextValue $ws.Range("D34") "1.677"
$ws.Range("E34").Value = "  +11.21%  "

# Row 35 - Stellar
Set-TextValue $ws.Range("D35") "0.08239"
$ws.Range("E35").Value = "  +3.29%  "

# Row 36 - Aptos
Set-TextValue $ws.Range("D36") "12.06"
$ws.Range("E36").Value = "  +16.38%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D37") "5.257"
$ws.Range("E37").Value = "  +10.93%  "

# Row 38 - TrustWalletToken
Set-TextValue $ws.Range("D38") "1.274"
$ws.Range("E38").Value = "  +3.52%  "

# Row 39 - FraxShare
Set-TextValue $ws.Range("D39") "8.668"
$ws.Range("E39").Value = "  +16.82%  "

# Row 40 - Hedera
Set-TextValue $ws.Range("D40") "0.06119"
$ws.Range("E40").Value = "  +4.36%  "

# Row 41 - VeChain
Set-TextValue $ws.Range("D41") "0.02220"
$ws.Range("E41").Value = "  +8.67%  "

# Row 42 - Algorand
Set-TextValue $ws.Range("D42") "0.2023"
$ws.Range("E42").Value = "  +7.86%  "

# Row 43 - Frax
Set-TextValue $ws.Range("D43") "0.9917"
$ws.Range("E43").Value = "  +4.29%  "

# Row 44 - TheSandbox
Set-TextValue $ws.Range("D44") "0.5922"
$ws.Range("E44").Value = "  +11.83%  "

# Row 45 - PancakeSwap
Set-TextValue $ws.Range("D45") "3.840"
$ws.Range("E45").Value = "  +8.77%  "

# Row 46 - EnergySwap
Set-TextValue $ws.Range("D46") "13.19"
$ws.Range("E46").Value = "  +7.81%  "

# Row 47 - was Decentraland, now Quant
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D47") "128.09"
$ws.Range("E47").Value = "  +8.42%  "

# Row 48 - was Quant, now Decentraland
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D48") "0.5692"
$ws.Range("E48").Value = "  +9.56%  "

# Row 49 - NEARProtocol
Set-TextValue $ws.Range("D49") "1.971"
$ws.Range("E49").Value = "  +8.58%  "

# Row 50 - Cronos
Set-TextValue $ws.Range("D50") "0.06823"
$ws.Range("E50").Value = "  +5.33%  "

# Row 51 - Aave
Set-TextValue $ws.Range("D51") "73.81"
$ws.Range("E51").Value = "  +9.36%  "
